$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44386
$ws.Range("D3").Value = 44386
$ws.Range("D4").Value = 44224
$ws.Range("D5").Value = 44224
$ws.Range("D6").Value = 44274
$ws.Range("D7").Value = 44274
$ws.Range("D8").Value = 44203
$ws.Range("D9").Value = 44203
$ws.Range("D10").Value = 44565
$ws.Range("O10").Value = 'Región Metropolitana'
$ws.Range("D11").Value = 44565
$ws.Range("O11").Value = 'Región Metropolitana'
$ws.Range("D12").Value = 44355
$ws.Range("D13").Value = 44355
$ws.Range("D14").Value = 44553
$ws.Range("D15").Value = 44553
$ws.Range("D16").Value = 44252
$ws.Range("D17").Value = 44252
$ws.Range("D18").Value = 44398
$ws.Range("O18").Value = 'Región de Ñuble'
$ws.Range("D19").Value = 44398
$ws.Range("O19").Value = 'Región de Ñuble'
$ws.Range("D20").Value = 44447
$ws.Range("D21").Value = 44447
$ws.Range("D22").Value = 44222
$ws.Range("D23").Value = 44222
$ws.Range("D24").Value = 44469
$ws.Range("D25").Value = 44469
$ws.Range("D26").Value = 44362
$ws.Range("J26").Value = 200
$ws.Range("D27").Value = 44362
$ws.Range("J27").Value = 100
$ws.Range("D28").Value = 44383
$ws.Range("D29").Value = 44383
$ws.Range("D30").Value = 44358
$ws.Range("D31").Value = 44358
$ws.Range("D32").Value = 44425
$ws.Range("D33").Value = 44425
$ws.Range("D34").Value = 44365
$ws.Range("D35").Value = 44365
$ws.Range("D36").Value = 44217
$ws.Range("D37").Value = 44217
$ws.Range("D38").Value = 44344
$ws.Range("N38").Value = '$/docena de 1 kilo'
$ws.Range("D39").Value = 44344
$ws.Range("N39").Value = '$/docena de 1 kilo'
$ws.Range("D40").Value = 44285
$ws.Range("D41").Value = 44285
$ws.Range("D42").Value = 44336
$ws.Range("D43").Value = 44336
$ws.Range("D44").Value = 44166
$ws.Range("D45").Value = 44166
$ws.Range("D46").Value = 44442
$ws.Range("J46").Value = 300
$ws.Range("O46").Value = 'Región de Ñuble'
$ws.Range("D47").Value = 44442
$ws.Range("J47").Value = 150
$ws.Range("O47").Value = 'Región de Ñuble'
$ws.Range("D48").Value = 44435
$ws.Range("D49").Value = 44435
$ws.Range("D50").Value = 44231
$ws.Range("D51").Value = 44231
$ws.Range("D52").Value = 44349
$ws.Range("O52").Value = 'Región Metropolitana'
$ws.Range("D53").Value = 44349
$ws.Range("O53").Value = 'Región Metropolitana'
$ws.Range("D54").Value = 44187
$ws.Range("D55").Value = 44187
$ws.Range("D56").Value = 44391
$ws.Range("D57").Value = 44391
$ws.Range("D58").Value = 44433
$ws.Range("D59").Value = 44433
$ws.Range("D60").Value = 44237
$ws.Range("D61").Value = 44237
$ws.Range("D62").Value = 44400
$ws.Range("D63").Value = 44400
$ws.Range("D64").Value = 44350
$ws.Range("D65").Value = 44350
$ws.Range("D66").Value = 44453
$ws.Range("D67").Value = 44453
$ws.Range("D68").Value = 44475
$ws.Range("D69").Value = 44475
$ws.Range("D70").Value = 44523
$ws.Range("D71").Value = 44523
$ws.Range("D72").Value = 44292
$ws.Range("D73").Value = 44292
$ws.Range("D74").Value = 44168
$ws.Range("D75").Value = 44168
$ws.Range("D76").Value = 44299
$ws.Range("O76").Value = 'Región de Ñuble'
$ws.Range("D77").Value = 44299
$ws.Range("O77").Value = 'Región de Ñuble'
$ws.Range("D78").Value = 44160
$ws.Range("D79").Value = 44160
$ws.Range("D80").Value = 44308
$ws.Range("D81").Value = 44308
$ws.Range("D82").Value = 44320
$ws.Range("D83").Value = 44320
$ws.Range("D84").Value = 44306
$ws.Range("D85").Value = 44306
$ws.Range("D86").Value = 44316
$ws.Range("D87").Value = 44316
$ws.Range("D88").Value = 44460
$ws.Range("D89").Value = 44460
$ws.Range("D90").Value = 44272
$ws.Range("D91").Value = 44272
$ws.Range("D92").Value = 44313
$ws.Range("D93").Value = 44313
$ws.Range("D94").Value = 44334
$ws.Range("D95").Value = 44334
$ws.Range("D96").Value = 44405
$ws.Range("D97").Value = 44405
$ws.Range("D98").Value = 44280
$ws.Range("J98").Value = 200
$ws.Range("D99").Value = 44280
$ws.Range("J99").Value = 100
$ws.Range("D100").Value = 44330
$ws.Range("D101").Value = 44330
$ws.Range("D102").Value = 44239
$ws.Range("D103").Value = 44239
$ws.Range("D104").Value = 44476
$ws.Range("D105").Value = 44476
$ws.Range("D106").Value = 44250
$ws.Range("O106").Value = 'Región de Arica y Parinacota'
$ws.Range("D107").Value = 44250
$ws.Range("O107").Value = 'Región de Arica y Parinacota'
$ws.Range("D108").Value = 44488
$ws.Range("D109").Value = 44488
$ws.Range("D110").Value = 44341
$ws.Range("D111").Value = 44341
$ws.Range("D112").Value = 44278
$ws.Range("J112").Value = 300
$ws.Range("D113").Value = 44278
$ws.Range("J113").Value = 150
$ws.Range("D114").Value = 44194
$ws.Range("N114").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("D115").Value = 44194
$ws.Range("N115").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("D116").Value = 44490
$ws.Range("D117").Value = 44490
$ws.Range("D118").Value = 44525
$ws.Range("D119").Value = 44525
$ws.Range("D120").Value = 44327
$ws.Range("D121").Value = 44327
$ws.Range("D122").Value = 44512
$ws.Range("D123").Value = 44512
